# Refresh cryptos list values (prices / 1h volume %) pulled from coinranking.com
# Row 41/42 also swap rank order: ImmutableX now above Kaspa
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.364.61'
$ws.Range("E2").Value = '  +2.18%  '

$ws.Range("D3").Value = '3.585.91'
$ws.Range("E3").Value = '  +0.19%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''242.30'
$ws.Range("E5").Value = '  +2.36%  '

$ws.Range("E6").Value = '  +1.09%  '

$ws.Range("D7").Value = '''1.71'
$ws.Range("E7").Value = '  +15.86%  '

$ws.Range("D8").Value = '''0.426'
$ws.Range("E8").Value = '  +6.76%  '

$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  +4.53%  '

$ws.Range("D11").Value = '3.585.55'
$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = '''44.57'
$ws.Range("E12").Value = '  +5.04%  '

$ws.Range("E14").Value = '  -0.20%  '

$ws.Range("D15").Value = '4.252.08'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '97.195.07'
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("E17").Value = '  +3.18%  '

$ws.Range("D18").Value = '''8.67'
$ws.Range("E18").Value = '  +11.92%  '

$ws.Range("D19").Value = '3.583.67'
$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("D20").Value = '''12.72'
$ws.Range("E20").Value = '  +1.56%  '

$ws.Range("D21").Value = '''18.07'
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").Value = '''0.529'
$ws.Range("E22").Value = '  +10.54%  '

$ws.Range("D23").Value = '''3.50'
$ws.Range("E23").Value = '  +1.09%  '

$ws.Range("D24").Value = '''516.86'
$ws.Range("E24").Value = '  +1.59%  '

$ws.Range("E25").Value = '  +5.20%  '

$ws.Range("E26").Value = '  +0.42%  '

$ws.Range("D27").Value = '''101.47'
$ws.Range("E27").Value = '  +6.43%  '

$ws.Range("E28").Value = '  +3.26%  '

$ws.Range("D29").Value = '3.778.65'
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("D30").Value = '''0.163'
$ws.Range("E30").Value = '  +14.43%  '

$ws.Range("E31").Value = '  -0.38%  '

$ws.Range("D32").Value = '''11.91'
$ws.Range("E32").Value = '  +3.88%  '

$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  -0.33%  '

$ws.Range("E34").Value = '  +3.86%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").Value = '''31.77'
$ws.Range("E36").Value = '  -0.13%  '

$ws.Range("D37").Value = '''8.89'
$ws.Range("E37").Value = '  +3.97%  '

$ws.Range("D38").Value = '''619.97'
$ws.Range("E38").Value = '  +5.98%  '

$ws.Range("D39").Value = '''0.569'
$ws.Range("E39").Value = '  +1.59%  '

$ws.Range("D40").Value = '''1.65'
$ws.Range("E40").Value = '  -1.55%  '

$ws.Range("B41").Value = 'ImmutableX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D41").Value = '''1.95'
$ws.Range("E41").Value = '  +8.52%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.155'
$ws.Range("E42").Value = '  +2.66%  '

$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").Value = '''0.925'
$ws.Range("E44").Value = '  +2.16%  '

$ws.Range("D45").Value = '''6.01'
$ws.Range("E45").Value = '  +4.56%  '

$ws.Range("E46").Value = '  +5.88%  '

$ws.Range("E47").Value = '  -0.10%  '

$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("D49").Value = '''0.407'
$ws.Range("E49").Value = '  +29.39%  '

$ws.Range("D50").Value = '''8.52'
$ws.Range("E50").Value = '  +4.54%  '

$ws.Range("E51").Value = '  +8.20%  '
